# Fruta / hortaliza, semanal
# Insert two new weekly price rows (421, 422) into the daily-logic
# subset sheet, pushing all existing rows from 421 downward by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 421-422; everything currently at/after row 421
# (through the old last row 478) shifts down to 423-480.
$ws.Range("A421:A422").EntireRow.Insert()

# --- New row 421: Naranja Lane Late, Primera ---
$ws.Cells.Item(421, 1).Value = 9
$ws.Cells.Item(421, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(421, 3).Value = "Metropolitana"
$ws.Cells.Item(421, 4).Value = 44491
$ws.Cells.Item(421, 5).Value = 13
$ws.Cells.Item(421, 6).Value = "Fruta"
$ws.Cells.Item(421, 7).Value = 100102
$ws.Cells.Item(421, 8).Value = "Cítricos"
$ws.Cells.Item(421, 9).Value = 100102005
$ws.Cells.Item(421, 10).Value = "Naranja"
$ws.Cells.Item(421, 11).Value = "Lane Late"
$ws.Cells.Item(421, 12).Value = "Primera"
$ws.Cells.Item(421, 13).Value = 720
$ws.Cells.Item(421, 14).Value = 7500
$ws.Cells.Item(421, 15).Value = 8000
$ws.Cells.Item(421, 16).Value = 7743
$ws.Cells.Item(421, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(421, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(421, 19).Value = 430
$ws.Cells.Item(421, 20).Value = 18

# --- New row 422: Naranja Navel Late, Primera ---
$ws.Cells.Item(422, 1).Value = 9
$ws.Cells.Item(422, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(422, 3).Value = "Metropolitana"
$ws.Cells.Item(422, 4).Value = 44491
$ws.Cells.Item(422, 5).Value = 13
$ws.Cells.Item(422, 6).Value = "Fruta"
$ws.Cells.Item(422, 7).Value = 100102
$ws.Cells.Item(422, 8).Value = "Cítricos"
$ws.Cells.Item(422, 9).Value = 100102005
$ws.Cells.Item(422, 10).Value = "Naranja"
$ws.Cells.Item(422, 11).Value = "Navel Late"
$ws.Cells.Item(422, 12).Value = "Primera"
$ws.Cells.Item(422, 13).Value = 33
$ws.Cells.Item(422, 14).Value = 150000
$ws.Cells.Item(422, 15).Value = 160000
$ws.Cells.Item(422, 16).Value = 155455
$ws.Cells.Item(422, 17).Value = "$/bins (400 kilos)"
$ws.Cells.Item(422, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(422, 19).Value = 389
$ws.Cells.Item(422, 20).Value = 400
